# Fruta / hortaliza, semanal
# Inserts two new weekly price rows for "Frutilla" (Terminal Hortofrutícola
# Agro Chillán) just above the existing row 228, pushing the subsequent
# rows (228-245) down to (230-247). The newly inserted rows hold a new
# "Especial" / "Primera" quality pair dated 44585 (2022-01-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 228; everything at/after row 228 shifts
# down by two rows (228->230 ... 245->247).
$ws.Rows("228:229").Insert()

# ---- Row 228 : Calidad "Especial" ----
$ws.Cells.Item(228, 1).Value = 7
$ws.Cells.Item(228, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(228, 3).Value = 'Ñuble'
$ws.Cells.Item(228, 4).Value = 44585
$ws.Cells.Item(228, 5).Value = 16
$ws.Cells.Item(228, 6).Value = 'Fruta'
$ws.Cells.Item(228, 7).Value = 100101
$ws.Cells.Item(228, 8).Value = 'Berries'
$ws.Cells.Item(228, 9).Value = 100112025
$ws.Cells.Item(228, 10).Value = 'Frutilla'
$ws.Cells.Item(228, 11).Value = 'Sin especificar'
$ws.Cells.Item(228, 12).Value = 'Especial'
$ws.Cells.Item(228, 13).Value = 60
$ws.Cells.Item(228, 14).Value = 7000
$ws.Cells.Item(228, 15).Value = 7000
$ws.Cells.Item(228, 16).Value = 7000
$ws.Cells.Item(228, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(228, 18).Value = 'Provincia de Diguillín'
$ws.Cells.Item(228, 19).Value = 1000
$ws.Cells.Item(228, 20).Value = 7

# ---- Row 229 : Calidad "Primera" ----
$ws.Cells.Item(229, 1).Value = 7
$ws.Cells.Item(229, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(229, 3).Value = 'Ñuble'
$ws.Cells.Item(229, 4).Value = 44585
$ws.Cells.Item(229, 5).Value = 16
$ws.Cells.Item(229, 6).Value = 'Fruta'
$ws.Cells.Item(229, 7).Value = 100101
$ws.Cells.Item(229, 8).Value = 'Berries'
$ws.Cells.Item(229, 9).Value = 100112025
$ws.Cells.Item(229, 10).Value = 'Frutilla'
$ws.Cells.Item(229, 11).Value = 'Sin especificar'
$ws.Cells.Item(229, 12).Value = 'Primera'
$ws.Cells.Item(229, 13).Value = 120
$ws.Cells.Item(229, 14).Value = 6000
$ws.Cells.Item(229, 15).Value = 6500
$ws.Cells.Item(229, 16).Value = 6250
$ws.Cells.Item(229, 17).Value = '$/caja 7 kilos'
$ws.Cells.Item(229, 18).Value = 'Provincia de Diguillín'
$ws.Cells.Item(229, 19).Value = 893
$ws.Cells.Item(229, 20).Value = 7

# Give column D of the new rows the same date number format the rest of
# the column uses.
$ws.Range("D228:D229").NumberFormat = $ws.Range("D230").NumberFormat
